$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column B into column C for rows 1-5 and 7 (row 6 intentionally skipped)
foreach ($r in 1,2,3,4,5,7) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 2).Text
}

# Update the active selection to match the final state
$ws.Range("L17").Select()
